# March 24 update 3
# Adds three new trailing columns to Sheet1:
#   M = renewd    (text,   constant "after" for every data row)
#   N = PlanID    (number, constant 20131419 for every data row)
#   O = iteration (number, constant 13 for every data row)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) ---
$ws.Range("M1").Value = "renewd"
$ws.Range("N1").Value = "PlanID"
$ws.Range("O1").Value = "iteration"

# Match the bold / centered / top-aligned / bordered style already used
# by the rest of the header row (B1:L1) by copying the format from L1.
$ws.Range("L1").Copy()
$ws.Range("M1:O1").PasteSpecial(-4122)   # xlPasteFormats

# --- Data rows (rows 2-31) ---
for ($row = 2; $row -le 31; $row++) {
    $ws.Cells.Item($row, 13).Value = "after"
    $ws.Cells.Item($row, 14).Value = 20131419
    $ws.Cells.Item($row, 15).Value = 13
}
